# Analyses des performances apres 7eme correction
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "LightHouse - Portable" ---
$ws1 = $wb.Worksheets.Item("LightHouse - Portable")
$ws1.Range("J4").Value = 83
$ws1.Range("J4").Select()

# --- Sheet 2: "LightHouse - Bureau" ---
$ws2 = $wb.Worksheets.Item("LightHouse - Bureau")
$ws2.Range("J4").Value = 90
$ws2.Range("J4").Interior.Color = $ws2.Range("E4").Interior.Color
$ws2.Range("J4").Select()

# --- Sheet 3: "GTmetrix - Bureau" ---
$ws3 = $wb.Worksheets.Item("GTmetrix - Bureau")
$ws3.Range("J4").Value = 100
$ws3.Range("J5").Value = 96
$ws3.Range("J5").Select()
